$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 87
$ws.Range("J2").Value = 334
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 97
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 54
$ws.Range("P2").Value = 0
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 38
$ws.Range("T2").Value = 64
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 603
$ws.Range("X2").Value = 609
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 6
